$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "69.379.91"
Set-TextValue $ws.Range("E2") "  +0.04%  "
Set-TextValue $ws.Range("D3") "3.667.12"
Set-TextValue $ws.Range("E3") "  -0.61%  "
Set-TextValue $ws.Range("E4") "  +0.01%  "
Set-TextValue $ws.Range("D5") "643.47"
Set-TextValue $ws.Range("E5") "  -5.53%  "
Set-TextValue $ws.Range("E6") "  -0.59%  "
Set-TextValue $ws.Range("E7") "  -0.01%  "
Set-TextValue $ws.Range("E8") "  +0.26%  "
Set-TextValue $ws.Range("E9") "  -1.30%  "
Set-TextValue $ws.Range("D10") "7.02"
Set-TextValue $ws.Range("E10") "  -1.46%  "
Set-TextValue $ws.Range("D11") "0.437"
Set-TextValue $ws.Range("E11") "  -0.30%  "
Set-TextValue $ws.Range("D12") "0.0000228"
Set-TextValue $ws.Range("E12") "  -1.64%  "
Set-TextValue $ws.Range("D13") "4.286.79"
Set-TextValue $ws.Range("E13") "  -0.58%  "
Set-TextValue $ws.Range("D14") "32.20"
Set-TextValue $ws.Range("E14") "  -0.91%  "
Set-TextValue $ws.Range("D15") "3.672.95"
Set-TextValue $ws.Range("E15") "  -0.54%  "
Set-TextValue $ws.Range("D16") "69.377.46"
Set-TextValue $ws.Range("E16") "  +0.03%  "
Set-TextValue $ws.Range("E17") "  +1.24%  "
Set-TextValue $ws.Range("D18") "15.91"
Set-TextValue $ws.Range("E18") "  -0.76%  "
Set-TextValue $ws.Range("D19") "6.39"
Set-TextValue $ws.Range("E19") "  -0.67%  "
Set-TextValue $ws.Range("D20") "463.91"
Set-TextValue $ws.Range("E20") "  -0.98%  "
Set-TextValue $ws.Range("D21") "9.68"
Set-TextValue $ws.Range("E21") "  -3.17%  "
Set-TextValue $ws.Range("D22") "0.640"
Set-TextValue $ws.Range("E22") "  -1.94%  "
Set-TextValue $ws.Range("D23") "79.36"
Set-TextValue $ws.Range("D24") "3.814.57"
Set-TextValue $ws.Range("E24") "  -0.57%  "
Set-TextValue $ws.Range("E25") "  -0.02%  "
Set-TextValue $ws.Range("E26") "  +0.40%  "
Set-TextValue $ws.Range("D27") "10.70"
Set-TextValue $ws.Range("E27") "  -2.45%  "
Set-TextValue $ws.Range("D28") "8.83"
Set-TextValue $ws.Range("E28") "  -3.50%  "
Set-TextValue $ws.Range("D29") "2.60"
Set-TextValue $ws.Range("E29") "  -3.11%  "
Set-TextValue $ws.Range("D30") "1.64"
Set-TextValue $ws.Range("E30") "  -6.40%  "
Set-TextValue $ws.Range("E31") "  +0.02%  "
Set-TextValue $ws.Range("D32") "1.97"
Set-TextValue $ws.Range("E32") "  -0.90%  "
Set-TextValue $ws.Range("D33") "26.46"
Set-TextValue $ws.Range("E33") "  -1.83%  "
Set-TextValue $ws.Range("D34") "6.39"
Set-TextValue $ws.Range("E34") "  -3.61%  "
Set-TextValue $ws.Range("D35") "3.658.10"
Set-TextValue $ws.Range("E35") "  -0.55%  "
Set-TextValue $ws.Range("E36") "  +1.58%  "
Set-TextValue $ws.Range("D37") "8.27"
Set-TextValue $ws.Range("E37") "  -0.39%  "
Set-TextValue $ws.Range("D39") "5.87"
Set-TextValue $ws.Range("E39") "  -5.90%  "
Set-TextValue $ws.Range("D40") "178.87"
Set-TextValue $ws.Range("E40") "  +4.18%  "
Set-TextValue $ws.Range("D41") "0.999"
Set-TextValue $ws.Range("E41") "  -0.10%  "
Set-TextValue $ws.Range("E42") "  -1.70%  "
Set-TextValue $ws.Range("E43") "  -4.28%  "
Set-TextValue $ws.Range("E44") "  -1.88%  "
Set-TextValue $ws.Range("D45") "46.57"
Set-TextValue $ws.Range("E45") "  -2.22%  "
Set-TextValue $ws.Range("D46") "2.68"
Set-TextValue $ws.Range("E46") "  -0.83%  "
Set-TextValue $ws.Range("D47") "26.87"
Set-TextValue $ws.Range("E47") "  -5.33%  "
Set-TextValue $ws.Range("E48") "  -3.41%  "
Set-TextValue $ws.Range("D49") "7.75"
Set-TextValue $ws.Range("E49") "  -0.62%  "
Set-TextValue $ws.Range("D50") "0.000261"
Set-TextValue $ws.Range("E50") "  -5.50%  "
Set-TextValue $ws.Range("E51") "  -6.55%  "
